$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.791.65'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.909.47'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '591.91'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.04%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '145.44'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.67%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '6.88'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').Value = '3.390.86'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '60.820.19'
$ws.Range('E16').Value = '  -0.09%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '6.66'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').Value = '2.910.72'
$ws.Range('E18').Value = '  -0.08%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '429.25'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '13.30'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  -1.21%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '81.36'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.32%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '10.94'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('E25').Value = '  -0.74%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.87'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.00%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.28'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +5.25%  '
$ws.Range('E29').Value = '  -0.02%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.60'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.56%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.02'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -3.22%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '26.46'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('D34').Value = '0.0₃0848'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.97'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.122'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.64%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '8.50'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('E41').Value = '  -2.32%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '39.85'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -4.34%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '373.71'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('D45').Value = '2.697.13'
$ws.Range('E45').Value = '  +0.88%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '132.16'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('E47').Value = '  -0.10%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '23.70'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -4.62%  '
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('E51').Value = '  +0.77%  '
